$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I1 = "I0" and J1 = "IF", matching the style of the existing
# header cell H1 (bold, centered, thin border).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-39: add values for columns I (I0) and J (IF)
$data = @{
  2  = @(4, 9)
  3  = @(7, 7)
  4  = @(5, 5)
  5  = @(9, 9)
  6  = @(4, 8)
  7  = @(1, 4)
  8  = @(1, 5)
  9  = @(1, 6)
  10 = @(1, 3)
  11 = @(1, 3)
  12 = @(1, 7)
  13 = @(1, 6)
  14 = @(1, 7)
  15 = @(1, 6)
  16 = @(1, 3)
  17 = @(1, 8)
  18 = @(1, 7)
  19 = @(1, 5)
  20 = @(1, 4)
  21 = @(1, 6)
  22 = @(1, 5)
  23 = @(1, 5)
  24 = @(1, 6)
  25 = @(1, 6)
  26 = @(1, 5)
  27 = @(1, 5)
  28 = @(1, 4)
  29 = @(1, 5)
  30 = @(1, 5)
  31 = @(1, 6)
  32 = @(1, 6)
  33 = @(1, 5)
  34 = @(1, 6)
  35 = @(1, 5)
  36 = @(1, 5)
  37 = @(5, 7)
  38 = @(1, 2)
  39 = @(1, 2)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Cells.Item($row, 9).Value = $vals[0]
  $ws.Cells.Item($row, 10).Value = $vals[1]
}
